$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the header formatting/style
# from the neighboring header cell (G1) so it matches the rest of the
# header row (bold, bordered, centered).
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the data values for the new Save column (plain/default style,
# like the rest of the data cells).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
